# Add two new stat columns, "I0" (col I) and "IF" (col J), to the sheet,
# matching the header style already used by the other header cells (e.g. H1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, centered, bordered) from H1
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:I70 and J2:J70
$iVals = @(7,8,7,8,8,10,8,7,8,8,7,7,7,6,7,7,6,10,8,7,7,8,7,7,8,7,7,9,7,7,7,6,5,5,7,12,7,7,6,8,8,7,9,8,8,8,7,5,5,7,7,6,7,6,8,9,8,9,7,7,8,7,8,7,5,6,5,5,8)
$jVals = @(7,8,7,8,8,10,8,7,8,8,7,7,7,6,7,7,6,10,8,7,7,8,7,7,8,8,8,9,8,7,8,7,6,6,8,12,8,7,6,8,8,7,9,8,8,8,8,6,5,7,7,7,7,7,9,9,8,9,7,7,8,8,8,7,5,6,5,5,8)

for ($idx = 0; $idx -lt $iVals.Count; $idx++) {
    $row = 2 + $idx
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
